$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Feuille1")
$ws2 = $wb.Worksheets.Item("Feuille2")

# --- Feuille1 ---
# the shared string used by A1 changes from "simple" to "multi"
$ws1.Range("A1").Value = "multi"

# --- Feuille2 ---
# new content: B1 (present, blank, default style) and B2 ("test")
$ws2.Range("B1").Font.Name = "Arial"
$ws2.Range("B2").Value = "test"

# update the active cell / selection on each sheet
$ws1.Range("B1").Select()
$ws2.Range("B2").Select()

# Feuille1 remains the active (tab-selected) sheet
$ws1.Activate()
